# 31.10.2020 MC Sales Details
# Fill in the retailer table (RetailerID / RetailerName / RetailerAddress)
# on the active sheet with the new batch of inactive-retailer records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 and row 4 were keyed in "Name, ID, Address" order; the rest in the
# natural "ID, Name, Address" column order (this ordering drives the
# resulting shared-strings table layout).
$ws.Cells.Item(2, 2).Value = "Azer Telecom"
$ws.Cells.Item(2, 1).Value = "RET-08809"
$ws.Cells.Item(2, 3).Value = "Arani BazarBaghaRajshahi "

$ws.Cells.Item(3, 1).Value = "RET-08821"
$ws.Cells.Item(3, 2).Value = "Lalon Telecom"
$ws.Cells.Item(3, 3).Value = "School MarketPuthiaRajshahi"

$ws.Cells.Item(4, 2).Value = "Ahona Telecom"
$ws.Cells.Item(4, 1).Value = "RET-20749"
$ws.Cells.Item(4, 3).Value = "Yusufpur Bazar Charghat Rajshahi"

$ws.Cells.Item(5, 1).Value = "RET-21074"
$ws.Cells.Item(5, 2).Value = "Piku Telecom"
$ws.Cells.Item(5, 3).Value = "Arani Bazar Bagha Rajshahi "

$ws.Cells.Item(6, 1).Value = "RET-21075"
$ws.Cells.Item(6, 2).Value = "Samsul Pharmacy "
$ws.Cells.Item(6, 3).Value = "Durduria Bazar Lalpur Natore "

$ws.Cells.Item(7, 1).Value = "RET-21146"
$ws.Cells.Item(7, 2).Value = "Rizia Variety Store "
$ws.Cells.Item(7, 3).Value = "Khanpur Bazar Bagha Rajshahi "

$ws.Cells.Item(8, 1).Value = "RET-29194"
$ws.Cells.Item(8, 2).Value = "Sampa Telecom"
$ws.Cells.Item(8, 3).Value = "Mirganj Bazar Bagha Rajshahi "

$ws.Cells.Item(9, 1).Value = "RET-29196"
$ws.Cells.Item(9, 2).Value = "Majumdar Electronics"
$ws.Cells.Item(9, 3).Value = "Lalpur Bazar Lalpur Natore "

$ws.Cells.Item(10, 1).Value = "RET-34481"
$ws.Cells.Item(10, 2).Value = "Liza Telecom"
$ws.Cells.Item(10, 3).Value = "Somospara Singra."

# Row 5 (RetailerID) and row 10 (RetailerAddress) had their formatting
# switched to the plain bordered (non-wrapping) style - reuse the exact
# existing style (same one already applied to A2) instead of creating a
# brand-new one.
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("L10").Select()
